$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.618
$ws.Range("E2").Value = 1.349
$ws.Range("F2").Value = 1.298
$ws.Range("G2").Value = 0.28
$ws.Range("H2").Value = 1.949
$ws.Range("D3").Value = 0.618
$ws.Range("E3").Value = 1.141
$ws.Range("F3").Value = 1.328
$ws.Range("G3").Value = 0.409
$ws.Range("H3").Value = 1.949
$ws.Range("D4").Value = 0.606
$ws.Range("E4").Value = 1.158
$ws.Range("F4").Value = 1.119
$ws.Range("G4").Value = 0.187
$ws.Range("H4").Value = 1.645
$ws.Range("D5").Value = 0.606
$ws.Range("E5").Value = 1.176
$ws.Range("F5").Value = 1.151
$ws.Range("G5").Value = 0.249
$ws.Range("H5").Value = 1.848
$ws.Range("D6").Value = 0.618
$ws.Range("E6").Value = 1.141
$ws.Range("F6").Value = 1.328
$ws.Range("G6").Value = 0.316
$ws.Range("H6").Value = 1.949
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 0.606
$ws.Range("E7").Value = 0.805
$ws.Range("F7").Value = 0.798
$ws.Range("G7").Value = 0.39
$ws.Range("H7").Value = 1.528
$ws.Range("D8").Value = 0.606
$ws.Range("E8").Value = 1.176
$ws.Range("F8").Value = 1.151
$ws.Range("G8").Value = 0.195
$ws.Range("H8").Value = 1.848
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 0.613
$ws.Range("E9").Value = 0.984
$ws.Range("F9").Value = 1.037
$ws.Range("G9").Value = 0.262
$ws.Range("H9").Value = 1.609
$ws.Range("D10").Value = 0.618
$ws.Range("E10").Value = 1.141
$ws.Range("F10").Value = 1.328
$ws.Range("G10").Value = 0.257
$ws.Range("H10").Value = 1.949
$ws.Range("C11").Value = 11
$ws.Range("D11").Value = 0.606
$ws.Range("E11").Value = 0.805
$ws.Range("F11").Value = 0.798
$ws.Range("G11").Value = 0.369
$ws.Range("H11").Value = 1.528
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 0.606
$ws.Range("E12").Value = 0.966
$ws.Range("F12").Value = 0.893
$ws.Range("G12").Value = 0.22
$ws.Range("H12").Value = 1.708
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 0.613
$ws.Range("E13").Value = 0.984
$ws.Range("F13").Value = 1.037
$ws.Range("G13").Value = 0.256
$ws.Range("H13").Value = 1.609
$ws.Range("D14").Value = 0.762
$ws.Range("E14").Value = 0.75
$ws.Range("F14").Value = 0.75
$ws.Range("G14").Value = 0.47
$ws.Range("H14").Value = 2.074
$ws.Range("D15").Value = 0.618
$ws.Range("E15").Value = 1.141
$ws.Range("F15").Value = 1.328
$ws.Range("G15").Value = 0.257
$ws.Range("H15").Value = 1.949
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 0.606
$ws.Range("E16").Value = 0.805
$ws.Range("F16").Value = 0.798
$ws.Range("G16").Value = 0.232
$ws.Range("H16").Value = 1.528
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 0.606
$ws.Range("E17").Value = 0.966
$ws.Range("F17").Value = 0.893
$ws.Range("G17").Value = 0.22
$ws.Range("H17").Value = 1.708
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 0.613
$ws.Range("E18").Value = 0.984
$ws.Range("F18").Value = 1.037
$ws.Range("G18").Value = 0.235
$ws.Range("H18").Value = 1.609
$ws.Range("D19").Value = 0.762
$ws.Range("E19").Value = 0.75
$ws.Range("F19").Value = 0.75
$ws.Range("G19").Value = 0.47
$ws.Range("H19").Value = 2.074
$ws.Range("D20").Value = 0.618
$ws.Range("E20").Value = 0.403
$ws.Range("F20").Value = 0.419
$ws.Range("G20").Value = 0.656
$ws.Range("H20").Value = 1.533
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 1.021
$ws.Range("E21").Value = 1.405
$ws.Range("F21").Value = 1.491
$ws.Range("G21").Value = 0.019
$ws.Range("H21").Value = 2.356
$ws.Range("C22").Value = 11
$ws.Range("D22").Value = 0.606
$ws.Range("E22").Value = 0.805
$ws.Range("F22").Value = 0.798
$ws.Range("G22").Value = 0.217
$ws.Range("H22").Value = 1.528
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 0.606
$ws.Range("E23").Value = 0.966
$ws.Range("F23").Value = 0.893
$ws.Range("G23").Value = 0.193
$ws.Range("H23").Value = 1.708
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 0.613
$ws.Range("E24").Value = 0.805
$ws.Range("F24").Value = 0.735
$ws.Range("G24").Value = 0.295
$ws.Range("H24").Value = 1.488
$ws.Range("D25").Value = 0.896
$ws.Range("E25").Value = 0.685
$ws.Range("F25").Value = 0.685
$ws.Range("G25").Value = 0.436
$ws.Range("H25").Value = 1.784
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 0.762
$ws.Range("E26").Value = 0.75
$ws.Range("F26").Value = 0.75
$ws.Range("G26").Value = 0.47
$ws.Range("H26").Value = 2.074
$ws.Range("D27").Value = 0.618
$ws.Range("E27").Value = 0.403
$ws.Range("F27").Value = 0.419
$ws.Range("G27").Value = 0.656
$ws.Range("H27").Value = 1.533
$ws.Range("D28").Value = 1.021
$ws.Range("E28").Value = 1.405
$ws.Range("F28").Value = 1.491
$ws.Range("G28").Value = -0.02
$ws.Range("H28").Value = 2.356
$ws.Range("C29").Value = 11
$ws.Range("D29").Value = 0.606
$ws.Range("E29").Value = 0.805
$ws.Range("F29").Value = 0.798
$ws.Range("G29").Value = 0.217
$ws.Range("H29").Value = 1.528
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 0.606
$ws.Range("E30").Value = 0.966
$ws.Range("F30").Value = 0.893
$ws.Range("G30").Value = 0.193
$ws.Range("H30").Value = 1.708
$ws.Range("C31").Value = 6
$ws.Range("D31").Value = 0.613
$ws.Range("E31").Value = 0.805
$ws.Range("F31").Value = 0.735
$ws.Range("G31").Value = 0.295
$ws.Range("H31").Value = 1.488
$ws.Range("D32").Value = 0.896
$ws.Range("E32").Value = 0.685
$ws.Range("F32").Value = 0.685
$ws.Range("G32").Value = 0.436
$ws.Range("H32").Value = 1.784
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = 0.762
$ws.Range("E33").Value = 0.75
$ws.Range("F33").Value = 0.75
$ws.Range("G33").Value = 0.47
$ws.Range("H33").Value = 2.074
$ws.Range("D34").Value = 0.618
$ws.Range("E34").Value = 0.403
$ws.Range("F34").Value = 0.419
$ws.Range("G34").Value = 0.653
$ws.Range("H34").Value = 1.533
$ws.Range("D35").Value = 1.021
$ws.Range("E35").Value = 1.122
$ws.Range("F35").Value = 1.122
$ws.Range("G35").Value = 0.125
$ws.Range("H35").Value = 2.157
$ws.Range("D36").Value = 1.443
$ws.Range("H36").Value = 2.645
$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 0.606
$ws.Range("E37").Value = 0.652
$ws.Range("F37").Value = 0.657
$ws.Range("G37").Value = 0.324
$ws.Range("H37").Value = 1.502
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 0.606
$ws.Range("E38").Value = 0.966
$ws.Range("F38").Value = 0.893
$ws.Range("G38").Value = 0.182
$ws.Range("H38").Value = 1.708
$ws.Range("C39").Value = 6
$ws.Range("D39").Value = 0.613
$ws.Range("E39").Value = 0.805
$ws.Range("F39").Value = 0.735
$ws.Range("G39").Value = 0.263
$ws.Range("H39").Value = 1.488
$ws.Range("D40").Value = 0.896
$ws.Range("E40").Value = 0.685
$ws.Range("F40").Value = 0.685
$ws.Range("G40").Value = 0.436
$ws.Range("H40").Value = 1.784
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = 0.762
$ws.Range("E41").Value = 0.75
$ws.Range("F41").Value = 0.75
$ws.Range("G41").Value = 0.47
$ws.Range("H41").Value = 2.074
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = 0.715
$ws.Range("E42").Value = 0.932
$ws.Range("F42").Value = 0.932
$ws.Range("G42").Value = 0.066
$ws.Range("H42").Value = 1.379
$ws.Range("D43").Value = 0.618
$ws.Range("E43").Value = 0.403
$ws.Range("F43").Value = 0.419
$ws.Range("G43").Value = 0.644
$ws.Range("H43").Value = 1.533
$ws.Range("D44").Value = 1.021
$ws.Range("E44").Value = 1.122
$ws.Range("F44").Value = 1.122
$ws.Range("G44").Value = 0.125
$ws.Range("H44").Value = 2.157
$ws.Range("D45").Value = 1.443
$ws.Range("H45").Value = 2.645
$ws.Range("C46").Value = 9
$ws.Range("D46").Value = 0.606
$ws.Range("E46").Value = 0.652
$ws.Range("F46").Value = 0.657
$ws.Range("G46").Value = 0.314
$ws.Range("H46").Value = 1.502
$ws.Range("C47").Value = 5
$ws.Range("D47").Value = 0.606
$ws.Range("E47").Value = 0.778
$ws.Range("F47").Value = 0.791
$ws.Range("G47").Value = 0.244
$ws.Range("H47").Value = 1.564
$ws.Range("C48").Value = 6
$ws.Range("D48").Value = 0.613
$ws.Range("E48").Value = 0.805
$ws.Range("F48").Value = 0.735
$ws.Range("G48").Value = 0.263
$ws.Range("H48").Value = 1.488
$ws.Range("D49").Value = 0.896
$ws.Range("E49").Value = 0.685
$ws.Range("F49").Value = 0.685
$ws.Range("G49").Value = 0.436
$ws.Range("H49").Value = 1.784
$ws.Range("C50").Value = 2
$ws.Range("D50").Value = 0.762
$ws.Range("E50").Value = 0.75
$ws.Range("F50").Value = 0.75
$ws.Range("G50").Value = 0.427
$ws.Range("H50").Value = 2.074
$ws.Range("C51").Value = 2
$ws.Range("D51").Value = 0.715
$ws.Range("E51").Value = 0.932
$ws.Range("F51").Value = 0.932
$ws.Range("G51").Value = 0.065
$ws.Range("H51").Value = 1.379
$ws.Range("C52").Value = 4
$ws.Range("D52").Value = 0.618
$ws.Range("E52").Value = 0.403
$ws.Range("F52").Value = 0.419
$ws.Range("G52").Value = 0.644
$ws.Range("H52").Value = 1.533
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 1.15
$ws.Range("E53").Value = "#NUM!"
$ws.Range("F53").ClearContents()
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 2.231
$ws.Range("D54").Value = 1.021
$ws.Range("E54").Value = 1.122
$ws.Range("F54").Value = 1.122
$ws.Range("G54").Value = 0.125
$ws.Range("H54").Value = 2.157
$ws.Range("D55").Value = 1.443
$ws.Range("H55").Value = 2.645
